# Applies the Odd_* value corrections from the 2025-03-31 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.63
$ws.Range("I2").Value = 2.9
$ws.Range("J2").Value = 1.08
$ws.Range("K2").Value = 8
$ws.Range("U2").Value = 12
$ws.Range("W2").Value = 26
$ws.Range("Y2").Value = 34
$ws.Range("AE2").Value = 8.5
$ws.Range("AF2").Value = 13
$ws.Range("AH2").Value = 29
$ws.Range("AI2").Value = 23

# Row 4
$ws.Range("N4").Value = 1.86
$ws.Range("O4").Value = 2.04

# Row 5
$ws.Range("L5").Value = 1.8
$ws.Range("M5").Value = 1.91
$ws.Range("N5").Value = 3.6
$ws.Range("O5").Value = 1.29
$ws.Range("P5").Value = 1.8
$ws.Range("Q5").Value = 2

# Row 6
$ws.Range("N6").Value = 2.2
$ws.Range("O6").Value = 1.65

# Row 7
$ws.Range("J7").Value = 1.1
$ws.Range("K7").Value = 7

# Row 8
$ws.Range("G8").Value = 1.44
$ws.Range("H8").Value = 4.5
$ws.Range("J8").Value = 1.04
$ws.Range("K8").Value = 13
$ws.Range("L8").Value = 1.22
$ws.Range("M8").Value = 4
$ws.Range("N8").Value = 1.73
$ws.Range("O8").Value = 2.08
$ws.Range("P8").Value = 1.33
$ws.Range("Q8").Value = 3.25
$ws.Range("R8").Value = 1.91
$ws.Range("S8").Value = 1.8
$ws.Range("T8").Value = 7
$ws.Range("X8").Value = 12
$ws.Range("Y8").Value = 26
$ws.Range("Z8").Value = 12
$ws.Range("AA8").Value = 8.5
$ws.Range("AG8").Value = 19
$ws.Range("AH8").Value = 67

# Row 10
$ws.Range("G10").Value = 4
$ws.Range("I10").Value = 1.8
$ws.Range("K10").Value = 21
$ws.Range("R10").Value = 1.44
$ws.Range("S10").Value = 2.63
$ws.Range("U10").Value = 23
$ws.Range("V10").Value = 13
$ws.Range("X10").Value = 26
$ws.Range("Y10").Value = 26
$ws.Range("Z10").Value = 21
$ws.Range("AF10").Value = 12
$ws.Range("AH10").Value = 17
$ws.Range("AI10").Value = 13

# Row 11
$ws.Range("G11").Value = 1.57
$ws.Range("H11").Value = 4.5
$ws.Range("I11").Value = 5
$ws.Range("K11").Value = 17
$ws.Range("R11").Value = 1.67
$ws.Range("S11").Value = 2.1
$ws.Range("W11").Value = 12
$ws.Range("AB11").Value = 15
$ws.Range("AI11").Value = 41
$ws.Range("AJ11").Value = 41

# Row 18
$ws.Range("I18").Value = 3.8
$ws.Range("Y18").Value = 21
$ws.Range("AA18").Value = 7.5
$ws.Range("AE18").Value = 15

# Row 21
$ws.Range("L21").Value = 1.29
$ws.Range("M21").Value = 3.5
$ws.Range("N21").Value = 1.95
$ws.Range("O21").Value = 1.9

# Row 22
$ws.Range("K22").Value = 9

# Row 28
$ws.Range("G28").Value = 2.55
$ws.Range("I28").Value = 2.7
$ws.Range("J28").Value = 1.08
$ws.Range("K28").Value = 8
$ws.Range("T28").Value = 7
$ws.Range("Z28").Value = 8
$ws.Range("AB28").Value = 17
$ws.Range("AD28").Value = 401
$ws.Range("AF28").Value = 13
$ws.Range("AG28").Value = 11
$ws.Range("AH28").Value = 29

# Row 30
$ws.Range("G30").Value = 3.15
$ws.Range("H30").Value = 3.25
$ws.Range("I30").Value = 2.2
$ws.Range("K30").Value = 7.2
$ws.Range("L30").Value = 1.3
$ws.Range("N30").Value = 1.88
$ws.Range("O30").Value = 1.83
$ws.Range("Q30").Value = 2.65
$ws.Range("T30").Value = 9.5
$ws.Range("U30").Value = 16.5
$ws.Range("V30").Value = 11
$ws.Range("W30").Value = 40
$ws.Range("X30").Value = 28
$ws.Range("Y30").Value = 35
$ws.Range("Z30").Value = 7.2
$ws.Range("AF30").Value = 11.25
$ws.Range("AG30").Value = 8.75
$ws.Range("AH30").Value = 22
$ws.Range("AI30").Value = 17
$ws.Range("AJ30").Value = 25

# Row 33
$ws.Range("G33").Value = 3.5
$ws.Range("H33").Value = 3.55
$ws.Range("J33").Value = 1.05
$ws.Range("K33").Value = 8
$ws.Range("L33").Value = 1.24
$ws.Range("M33").Value = 3.7
$ws.Range("N33").Value = 1.72
$ws.Range("O33").Value = 2.02
$ws.Range("P33").Value = 1.36
$ws.Range("Q33").Value = 2.87
$ws.Range("S33").Value = 2.15
$ws.Range("T33").Value = 12
$ws.Range("W33").Value = 50
$ws.Range("X33").Value = 28
$ws.Range("Y33").Value = 32
$ws.Range("Z33").Value = 8
$ws.Range("AA33").Value = 6.9
$ws.Range("AB33").Value = 12.5
$ws.Range("AD33").Value = 350
$ws.Range("AE33").Value = 8.5
$ws.Range("AG33").Value = 8.25
$ws.Range("AH33").Value = 18

# Row 36
$ws.Range("H36").Value = 3.1
$ws.Range("I36").Value = 5
$ws.Range("J36").Value = 1.1
$ws.Range("K36").Value = 7
$ws.Range("L36").Value = 1.44
$ws.Range("M36").Value = 2.63
$ws.Range("N36").Value = 2.35
$ws.Range("O36").Value = 1.57
$ws.Range("P36").Value = 1.5
$ws.Range("Q36").Value = 2.5
$ws.Range("R36").Value = 2.2
$ws.Range("S36").Value = 1.62
$ws.Range("T36").Value = 5.5
$ws.Range("U36").Value = 7.5
$ws.Range("X36").Value = 17
$ws.Range("Y36").Value = 34
$ws.Range("Z36").Value = 7
$ws.Range("AB36").Value = 19
$ws.Range("AC36").Value = 67
$ws.Range("AE36").Value = 11
$ws.Range("AJ36").Value = 51
